$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ministry Course Code and Level" column (G) is being split into two
# separate columns: "Ministry Course Code" and "Ministry Course Level".
# Insert a new column at H so everything from the old H onward shifts right.
$ws.Columns("H").Insert()

# Update the header row for the newly split columns.
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("H1").Value = "Ministry Course Level"

# Split the combined "ENST 12" values into a text course code (G) and a
# numeric course level (H) for every data row.
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 7).Value = "ENST"
    $ws.Cells.Item($r, 8).Value = 12
}

# Reflect the column split in the current selection.
$ws.Range("H1").Select() | Out-Null
